# ADD results from server
# Rewrite the header row (row 1) and the data row (row 2) on every
# year-sheet (2025..2050) with the new "gb"/"btes" columns and the
# updated numbers that came back from the server run.

$wb = $excel.ActiveWorkbook

# New header order shared by every sheet: "eb" stays in col A, "ites"
# stays in col O, but a new "gb" column is inserted right after "eb",
# "gt"/"dgt" are dropped, and a new "btes" column is inserted before
# "ites".
$headers = @("eb","gb","hp","st","wi","ieh","chp","ac","ab_ct","ab_hp","cp_ct","cp_hp","ttes","btes","ites")

# New row-2 numbers (A..O) per sheet, in the new column order above.
$valuesBySheet = @{
    1 = @(3906.399109145206,0,48353.76274462014,0,289724.0114301849,9433.134471502228,0,2534.277928792104,0,0,0,0,0,2367.37219622158,1995.762462679798)
    2 = @(6991.052031681918,0,197913.7502057619,0,289724.0114301849,16452.51445364119,0,8194.52068131253,0,0,0,0,0,7543.193583625169,6257.586732772244)
    3 = @(31236.29455387744,0,292247.2772138842,0,289724.0114301849,16595.10705160327,0,12131.91920790125,0,0,0,0,0,12888.04225687751,9263.466444480218)
    4 = @(31236.29455387744,0,292247.2772138842,0,289724.0114301849,16595.10705160327,0,12131.91920790125,0,0,0,0,0,14045.89200932069,9263.466444480218)
    5 = @(38906.8534480406,193.0947398408091,292247.2772138842,0,289724.0114301849,16595.10705160327,0,12131.91920790125,0,0,0,0,0,16879.89729726143,10096.02314047837)
    6 = @(38906.8534480406,193.0947398408091,292247.2772138842,0,289724.0114301849,16595.10705160327,0,12131.91920790125,0,0,0,0,0,14045.89200932069,9263.466444480218)
}

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }

    $rowValues = $valuesBySheet[$s]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $ws.Cells.Item(2, $i + 1).Value = $rowValues[$i]
    }
}
